$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 136
$ws.Range("I2").Value = 456
$ws.Range("J2").Value = 1719
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 505
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 291
$ws.Range("O2").Value = 2
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 27
$ws.Range("S2").Value = 166
$ws.Range("T2").Value = 276
$ws.Range("U2").Value = 23
$ws.Range("V2").Value = 2550
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2642
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 14
